# ran resolve and classify+summarise steps after changes to mapping file
$wb = $excel.ActiveWorkbook

# --- "Range Status" sheet ---
# All "Species (perc.)" (column C) values cleared; "Species (no.)" (column B)
# values reset to 0 for rows 2,4,5,6,7 (row 3 stays 0, only loses its C value).
$wsRange = $wb.Worksheets.Item("Range Status")
$wsRange.Range("B2").Value = 0
$wsRange.Range("C2").ClearContents()
$wsRange.Range("C3").ClearContents()
$wsRange.Range("B4").Value = 0
$wsRange.Range("C4").ClearContents()
$wsRange.Range("B5").Value = 0
$wsRange.Range("C5").ClearContents()
$wsRange.Range("B6").Value = 0
$wsRange.Range("C6").ClearContents()
$wsRange.Range("B7").Value = 0
$wsRange.Range("C7").ClearContents()

# --- "Priority Status" sheet ---
$wsPriority = $wb.Worksheets.Item("Priority Status")
$wsPriority.Range("B2").Value = 2
$wsPriority.Range("B3").Value = 44

# --- "Species qualification" sheet ---
$wsSpecies = $wb.Worksheets.Item("Species qualification")
$wsSpecies.Range("B5").Value = 0

# --- "High Priority break-up" sheet ---
$wsBreakup = $wb.Worksheets.Item("High Priority break-up")
$wsBreakup.Range("B2").Value = 2
$wsBreakup.Range("D2").Value = 2
